$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '70.699.21'
$ws.Range("E2").Value = '  -0.24%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.524.17'
$ws.Range("E3").Value = '  -1.05%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.10%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '610.98'
$ws.Range("E5").Value = '  -0.83%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '173.61'
$ws.Range("E6").Value = '  +0.79%  '

$ws.Range("B7").Value = 'LidoStakedEther'
$ws.Range("C7").Value = 'https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '3.518.42'
$ws.Range("E7").Value = '  -1.08%  '

$ws.Range("B8").Value = 'XRP'
$ws.Range("C8").Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.610'
$ws.Range("E8").Value = '  -1.31%  '

$ws.Range("E9").Value = '  -0.05%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.196'
$ws.Range("E10").Value = '  -0.49%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '7.31'
$ws.Range("E11").Value = '  +1.81%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.588'
$ws.Range("E12").Value = '  +0.28%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '46.49'
$ws.Range("E13").Value = '  -0.64%  '

$ws.Range("E14").Value = '  -0.64%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '4.091.96'
$ws.Range("E15").Value = '  -1.15%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '8.44'
$ws.Range("E16").Value = '  +0.48%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '615.83'
$ws.Range("E17").Value = '  -0.57%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.522.14'
$ws.Range("E18").Value = '  -0.33%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '70.720.72'
$ws.Range("E19").Value = '  -0.34%  '

$ws.Range("E20").Value = '  +1.24%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '17.78'
$ws.Range("E21").Value = '  +2.21%  '

$ws.Range("E22").Value = '  +0.15%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.00'
$ws.Range("E23").Value = '  -5.43%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '15.78'
$ws.Range("E24").Value = '  +0.31%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '98.05'
$ws.Range("E25").Value = '  +1.32%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.77'
$ws.Range("E26").Value = '  -1.50%  '

$ws.Range("E27").Value = '  +0.03%  '

$ws.Range("E28").Value = '  -0.32%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '33.70'
$ws.Range("E29").Value = '  +0.59%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '9.16'
$ws.Range("E30").Value = '  +0.78%  '

$ws.Range("E31").Value = '  -1.02%  '

$ws.Range("E32").Value = '  -4.07%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.30'
$ws.Range("E33").Value = '  -0.62%  '

$ws.Range("E34").Value = '  -2.15%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '616.02'
$ws.Range("E35").Value = '  +7.03%  '

$ws.Range("E36").Value = '  -0.95%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '10.86'
$ws.Range("E37").Value = '  -0.10%  '

$ws.Range("B38").Value = 'dogwifhat'
$ws.Range("C38").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.51'
$ws.Range("E38").Value = '  -2.94%  '

$ws.Range("B39").Value = 'VeChain'
$ws.Range("C39").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0477'
$ws.Range("E39").Value = '  +0.87%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '56.98'
$ws.Range("E40").Value = '  -1.02%  '

$ws.Range("E41").Value = '  +0.05%  '

$ws.Range("E42").Value = '  +1.53%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '3.375.57'
$ws.Range("E43").Value = '  +0.44%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0₃0737'
$ws.Range("E44").Value = '  +4.39%  '

$ws.Range("E45").Value = '  -2.41%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.95'
$ws.Range("E46").Value = '  -2.64%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '32.27'
$ws.Range("E47").Value = '  -2.58%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.58'
$ws.Range("E48").Value = '  -2.03%  '

$ws.Range("E49").Value = '  +0.28%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '133.64'
$ws.Range("E50").Value = '  -0.36%  '
